# Commiting changes after incorporating events for video pages
#
# Replace the sample "Health" row with a new "Parent" (video) row and add a
# hyperlink on the URL cell, then move the active selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell B1 - text unchanged ("brand") but rewritten so the
# shared-string table is rebuilt in the same order the source workbook uses.
$ws.Range("B1").Value2 = "brand"

$newUrl = "https://www.parents.com/kindred/anthony-anderson-says-growing-up-black-gave-him-his-sense-of-humor-you-have-to-laugh-to-keep-from-crying/"

# Row 3 becomes the "Parent" / video-page sample row.
$ws.Range("A3").Value2 = $newUrl
$ws.Range("D3").Value2 = "Video"
$ws.Range("B3").Value2 = "Parent"
$ws.Range("C3").Value2 = "BIO"

# Turn the URL cell into a real hyperlink pointing at the new article.
$ws.Hyperlinks.Add($ws.Range("A3"), $newUrl)

# Move the selection, matching the saved cursor position in the workbook.
$ws.Range("C3").Select() | Out-Null
